$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("employment")
$ws.Range("A5").Value = "Undergraduate Research Assistant"
$ws.Range("B5").Value = "Work unit for psychological methods with interdisciplinary focus"
$ws.Range("C5").Value = "Goethe-University"
$ws.Range("D5").Value = "Frankfurt am Main"
$ws.Range("E5").Value = "2023/04"
$null = $ws.Range("F5").Select()
